$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text-valued cells (coin names, links, percentages, multi-dot prices)
$ws.Range('D2').Value = '37.565.70'
$ws.Range('E2').Value = '  +1.94%  '
$ws.Range('D3').Value = '2.042.76'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('E5').Value = '  +5.06%  '
$ws.Range('E6').Value = '  -0.79%  '
$ws.Range('E8').Value = '  -5.72%  '
$ws.Range('E9').Value = '  +1.83%  '
$ws.Range('E10').Value = '  -0.50%  '
$ws.Range('E11').Value = '  -0.63%  '
$ws.Range('E12').Value = '  -1.02%  '
$ws.Range('E13').Value = '  +1.38%  '
$ws.Range('D14').Value = '2.346.07'
$ws.Range('E14').Value = '  +3.37%  '
$ws.Range('E15').Value = '  -2.56%  '
$ws.Range('E16').Value = '  -2.57%  '
$ws.Range('E17').Value = '  -0.84%  '
$ws.Range('D18').Value = '2.053.78'
$ws.Range('E18').Value = '  +3.82%  '
$ws.Range('D19').Value = '37.516.69'
$ws.Range('E19').Value = '  +2.01%  '
$ws.Range('E20').Value = '  +0.02%  '
$ws.Range('D21').Value = '0.0₃0857'
$ws.Range('E21').Value = '  -0.65%  '
$ws.Range('E22').Value = '  +1.81%  '
$ws.Range('E23').Value = '  -0.40%  '
$ws.Range('E24').Value = '  +7.71%  '
$ws.Range('E25').Value = '  -0.07%  '
$ws.Range('E26').Value = '  -0.81%  '
$ws.Range('E27').Value = '  -0.75%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('E29').Value = '  -4.41%  '
$ws.Range('E30').Value = '  +2.22%  '
$ws.Range('E31').Value = '  +0.30%  '
$ws.Range('E32').Value = '  -0.59%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('E33').Value = '  +7.71%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('E34').Value = '  -1.73%  '
$ws.Range('E35').Value = '  +0.47%  '
$ws.Range('E36').Value = '  +10.71%  '
$ws.Range('E37').Value = '  +3.94%  '
$ws.Range('E38').Value = '  -0.18%  '
$ws.Range('E39').Value = '  +2.45%  '
$ws.Range('E40').Value = '  -1.70%  '
$ws.Range('E41').Value = '  +4.06%  '
$ws.Range('E42').Value = '  -1.07%  '
$ws.Range('E43').Value = '  +2.90%  '
$ws.Range('E44').Value = '  +1.79%  '
$ws.Range('E45').Value = '  +1.32%  '
$ws.Range('D46').Value = '1.408.18'
$ws.Range('E46').Value = '  +2.92%  '
$ws.Range('E47').Value = '  +1.88%  '
$ws.Range('E48').Value = '  +1.88%  '
$ws.Range('E49').Value = '  +3.38%  '
$ws.Range('E50').Value = '  +2.38%  '
$ws.Range('E51').Value = '  +4.21%  '

# Numeric-looking price cells: force Text storage so they keep their exact
# string form (e.g. "5.39") instead of being parsed into a float, then reset
# the cell style back to Normal so no stray number-format style lingers.
$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '257.61'
$cell.Style = "Normal"
$cell = $ws.Range('D8')
$cell.NumberFormat = "@"
$cell.Value = '57.38'
$cell.Style = "Normal"
$cell = $ws.Range('D11')
$cell.NumberFormat = "@"
$cell.Value = '0.0798'
$cell.Style = "Normal"
$cell = $ws.Range('D15')
$cell.NumberFormat = "@"
$cell.Value = '0.823'
$cell.Style = "Normal"
$cell = $ws.Range('D16')
$cell.NumberFormat = "@"
$cell.Value = '21.42'
$cell.Style = "Normal"
$cell = $ws.Range('D17')
$cell.NumberFormat = "@"
$cell.Value = '5.39'
$cell.Style = "Normal"
$cell = $ws.Range('D22')
$cell.NumberFormat = "@"
$cell.Value = '5.24'
$cell.Style = "Normal"
$cell = $ws.Range('D23')
$cell.NumberFormat = "@"
$cell.Value = '229.19'
$cell.Style = "Normal"
$cell = $ws.Range('D27')
$cell.NumberFormat = "@"
$cell.Value = '9.19'
$cell.Style = "Normal"
$cell = $ws.Range('D28')
$cell.NumberFormat = "@"
$cell.Value = '163.38'
$cell.Style = "Normal"
$cell = $ws.Range('D29')
$cell.NumberFormat = "@"
$cell.Value = '0.139'
$cell.Style = "Normal"
$cell = $ws.Range('D30')
$cell.NumberFormat = "@"
$cell.Value = '19.91'
$cell.Style = "Normal"
$cell = $ws.Range('D33')
$cell.NumberFormat = "@"
$cell.Value = '0.0667'
$cell.Style = "Normal"
$cell = $ws.Range('D34')
$cell.NumberFormat = "@"
$cell.Value = '4.75'
$cell.Style = "Normal"
$cell = $ws.Range('D36')
$cell.NumberFormat = "@"
$cell.Value = '2.52'
$cell.Style = "Normal"
$cell = $ws.Range('D37')
$cell.NumberFormat = "@"
$cell.Value = '3.46'
$cell.Style = "Normal"
$cell = $ws.Range('D42')
$cell.NumberFormat = "@"
$cell.Value = '0.0966'
$cell.Style = "Normal"
$cell = $ws.Range('D45')
$cell.NumberFormat = "@"
$cell.Value = '16.32'
$cell.Style = "Normal"
$cell = $ws.Range('D47')
$cell.NumberFormat = "@"
$cell.Value = '91.41'
$cell.Style = "Normal"
$cell = $ws.Range('D48')
$cell.NumberFormat = "@"
$cell.Value = '1.06'
$cell.Style = "Normal"
$cell = $ws.Range('D49')
$cell.NumberFormat = "@"
$cell.Value = '7.48'
$cell.Style = "Normal"
$cell = $ws.Range('D51')
$cell.NumberFormat = "@"
$cell.Value = '2.02'
$cell.Style = "Normal"
